$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at U, shifting old U/V to V/W
$ws.Columns("U:U").Insert()

# Set header for the newly inserted column
$ws.Range("U1").Value = "MgCa Coretop modelled temperature"

# Update row 2 values
$ws.Range("R2").Value = 25.31
$ws.Range("S2").Value = -0.8899967108832421
$ws.Range("T2").Value = 0.4471032891167575
$ws.Range("U2").Value = 25.7765
$ws.Range("V2").Value = -1.356499999999997
$ws.Range("W2").Value = -0.01939999999999742
